$d = $word.ActiveDocument

# The paragraph "<div><id>p155r_1</id><head>..." currently stores the
# tag text "<id>p155r_1</id>" split across three runs with different
# character formatting:
#   run 1: "<id>"     (Courier New, color 7f6000, sz 18)
#   run 2: "p155r_1"  (default formatting, color 000000)
#   run 3: "</id>"    (Courier New, color 7f6000, sz 18)
#
# The edit merges them into a single run containing the full literal
# text "<id>p155r_1</id>", keeping the Courier New / 7f6000 / sz 18
# character formatting of the surrounding tag runs (i.e. the
# formatting of the first of the three runs, which Word's Find
# carries onto the replacement text).
#
# "<id>p155r_1</id>" is a unique string in the document (a similarly
# named "<id>fig_p155r_1</id>" elsewhere must stay untouched), so a
# plain literal Find/Replace unambiguously targets the right spot.

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.MatchCase = $true
$rng.Find.Execute("<id>p155r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p155r_1</id>", 2)
